# Rename the single worksheet from "updated" to "Tabelle1"
# (clean-up of input tables)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("updated")
$ws.Name = "Tabelle1"
